$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing indicator values (D2, D6) ---
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 2

# A new column (E) shows up carrying a 4-decimal numeric format; E6 stays
# empty but already formatted that way.
$ws.Range("E6").NumberFormat = "0.0000"
$ws.Columns.Item(5).ColumnWidth = 8.23

# --- New accounts: DVD and CAPT for Google/2014 ---
$ws.Range("A7").Value = "Google"
$ws.Range("B7").Value = 2014
$ws.Range("C7").Value = "DVD"
$ws.Range("D7").Value = 1

$ws.Range("A8").Value = "Google"
$ws.Range("B8").Value = 2014
$ws.Range("C8").Value = "CAPT"
$ws.Range("D8").Value = 0.5

# Leave selection where the author's cursor ended up after the edits
$ws.Range("D9").Select() | Out-Null
